# Updated index regressions with controls to include SBAC test scores
#
# The "Student Motivation" category index (and its blank separator row)
# is dropped from the multivariate regression table, "Student Support" is
# relabeled "Counseling Support", the summary note text is updated to
# reflect 3 (rather than 4) index regressors, and the N row is
# reformatted as centered whole numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Student Motivation" rows (data rows 15:16) together with the
# blank spacer row that followed them (row 17). Everything below shifts up.
$ws.Rows("15:17").Delete() | Out-Null

# Rename the "Student Support" category to "Counseling Support".
$ws.Range("A12").Value = "Counseling Support"

# Update the closing note to mention 3 (not 4) index variables and fix the
# "Reregressions" -> "Regressions" wording.
$ws.Range("A23").Value = "NOTE: Regressions are run using standardized z scores for all variables. These are multivariate regressions with all 3 index variables as regressors."

# The "N" row (now row 18) is reformatted with an integer number format and
# centered horizontal alignment.
$nRow = $ws.Range("B18:O18")
$nRow.NumberFormat = "0"
$nRow.HorizontalAlignment = -4108

# Leave the selection where the author left it after editing the sheet.
$ws.Range("A24").Select() | Out-Null
